# Expense tracker update:
#   - Row 2: "Food"/400 -> "Groceries"/40, with a refreshed timestamp
#   - Row 3 (new): "Spotify Subscription"/100, same timestamp, same date
#     number format as the Date column above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45965.22928240741

# Update existing row 2.
$ws.Range("A2").Value = "Groceries"
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = $newDate

# Give the new C3 cell the same date formatting as C2 before filling it in,
# so it picks up the existing date style instead of creating a new one.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the new row 3.
$ws.Range("A3").Value = "Spotify Subscription"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = $newDate
